# Apply updated Excess Mortality figures for a handful of rows
# (Czechia, Iceland, Estonia) following a recalculation of the underlying
# mortality numbers for 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: Female / Czechia -------------------------------------------------
$ws.Range("D7").Value = 22234
$ws.Range("P7").Value = 3291.6
$ws.Range("S7").Value = 1.9
$ws.Range("U7").Value = "3291.6 (" + [char]177 + "315.9)"
$ws.Range("V7").Value = "17.4% (" + [char]177 + "1.9%)"
$ws.Range("X7").Value = 64.2
$ws.Range("Y7").Value = 6.1
$ws.Range("Z7").Value = "64.2(" + [char]177 + "6.1)"

# --- Row 14: Female / Iceland -------------------------------------------------
$ws.Range("D14").Value = 339
$ws.Range("P14").Value = -0.4
$ws.Range("R14").Value = -0.1
$ws.Range("S14").Value = 4.6
$ws.Range("U14").Value = "-0.4 (" + [char]177 + "16.2)"
$ws.Range("V14").Value = "-0.1% (" + [char]177 + "4.6%)"
$ws.Range("X14").Value = -0.2
$ws.Range("Z14").Value = "-0.2(" + [char]177 + "9.5)"

# --- Row 36: Male / Czechia ---------------------------------------------------
$ws.Range("D36").Value = 37233
$ws.Range("P36").Value = 6018.4
$ws.Range("U36").Value = "6018.4 (" + [char]177 + "396.7)"

# --- Row 38: Male / Estonia ---------------------------------------------------
$ws.Range("D38").Value = 4272
$ws.Range("P38").Value = 60.2
$ws.Range("U38").Value = "60.2 (" + [char]177 + "110.4)"
$ws.Range("X38").Value = 9.9
$ws.Range("Z38").Value = "9.9(" + [char]177 + "18.1)"

# --- Row 43: Male / Iceland ----------------------------------------------------
$ws.Range("D43").Value = 539
$ws.Range("P43").Value = 56.2
$ws.Range("R43").Value = 11.6
$ws.Range("U43").Value = "56.2 (" + [char]177 + "18.1)"
$ws.Range("V43").Value = "11.6% (" + [char]177 + "4.0%)"
$ws.Range("X43").Value = 31
$ws.Range("Y43").Value = 9.9
$ws.Range("Z43").Value = "31.0(" + [char]177 + "9.9)"

# --- Row 65: Total / Czechia ---------------------------------------------------
$ws.Range("D65").Value = 59467
$ws.Range("P65").Value = 9310
$ws.Range("U65").Value = "9310.0 (" + [char]177 + "711.7)"

# --- Row 67: Total / Estonia ---------------------------------------------------
$ws.Range("D67").Value = 6722
$ws.Range("P67").Value = 50.4
$ws.Range("R67").Value = 0.8
$ws.Range("S67").Value = 2.9
$ws.Range("U67").Value = "50.4 (" + [char]177 + "191.8)"
$ws.Range("V67").Value = "0.8% (" + [char]177 + "2.9%)"
$ws.Range("X67").Value = 4
$ws.Range("Z67").Value = "4.0(" + [char]177 + "15.4)"

# --- Row 72: Total / Iceland ----------------------------------------------------
$ws.Range("D72").Value = 878
$ws.Range("P72").Value = 55.8
$ws.Range("R72").Value = 6.8
$ws.Range("U72").Value = "55.8 (" + [char]177 + "22.4)"
$ws.Range("V72").Value = "6.8% (" + [char]177 + "2.8%)"
$ws.Range("X72").Value = 15.9
$ws.Range("Y72").Value = 6.3
$ws.Range("Z72").Value = "15.9(" + [char]177 + "6.3)"
